# Apply the coin-price-list update described by the commit diff.
# D (Price) and E (Volume/1h %) columns hold text that LOOKS numeric
# (e.g. "308.60", "-1.38%") -- the source workbook stores them as plain
# text (inlineStr), so we prefix new values with a literal apostrophe to
# force Excel to keep them as text instead of parsing them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''308.60'
$ws.Range("E2").Value = '''-1.38%'

# Row 3
$ws.Range("D3").Value = '''36.94'
$ws.Range("E3").Value = '''-2.82%'

# Row 4
$ws.Range("D4").Value = '''5.137'
$ws.Range("E4").Value = '''-0.09%'

# Row 5
$ws.Range("D5").Value = '''0.07779'
$ws.Range("E5").Value = '''-1.86%'

# Row 6
$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").Value = '''8.308'
$ws.Range("E6").Value = '''0.53%'

# Row 7
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = '''1.869'
$ws.Range("E7").Value = '''-2.38%'

# Row 8
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").Value = '''2.983'
$ws.Range("E8").Value = '''-3.76%'

# Row 9
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '''0.9245'
$ws.Range("E9").Value = '''-0.21%'

# Row 10
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '''0.1113'
$ws.Range("E10").Value = '''-7.68%'

# Row 11
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1880'
$ws.Range("E11").Value = '''-1.90%'

# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.08843'
$ws.Range("E12").Value = '''-3.84%'

# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.03294'
$ws.Range("E13").Value = '''-1.57%'

# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.09577'
$ws.Range("E14").Value = '''-0.62%'

# Row 15
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001387'
$ws.Range("E15").Value = '''1.46%'

# Row 16
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.006021'
$ws.Range("E16").Value = '''1.65%'

# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.393'
$ws.Range("E17").Value = '''-4.07%'

# Row 18
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = '''4.401'
$ws.Range("E18").Value = '''-0.11%'

# Row 19
$ws.Range("E19").Value = '''0.06%'

# Row 20
$ws.Range("D20").Value = '''6.365'
$ws.Range("E20").Value = '''20.36%'

# Row 21
$ws.Range("E21").Value = '''0.57%'

# Row 22
$ws.Range("D22").Value = '''0.2375'
$ws.Range("E22").Value = '''-8.29%'

# Row 23
$ws.Range("D23").Value = '''0.04349'
$ws.Range("E23").Value = '''-0.39%'

# Row 24
$ws.Range("D24").Value = '''0.001203'
$ws.Range("E24").Value = '''-3.84%'

# Row 25
$ws.Range("D25").Value = '''0.004274'
$ws.Range("E25").Value = '''-0.54%'

# Row 26
$ws.Range("E26").Value = '''7.97%'

# Row 27
$ws.Range("D27").Value = '''0.0002906'

# Row 39
$ws.Range("D39").Value = '''0.02144'
$ws.Range("E39").Value = '''1.38%'

# Row 40
$ws.Range("D40").Value = '''0.04937'

# Row 41
$ws.Range("D41").Value = '''0.007571'
$ws.Range("E41").Value = '''-0.90%'

# Row 42
$ws.Range("E42").Value = '''-0.48%'

# Row 43
$ws.Range("D43").Value = '''0.008502'
$ws.Range("E43").Value = '''-6.75%'

# Row 44
$ws.Range("D44").Value = '''0.001991'
$ws.Range("E44").Value = '''-2.84%'

# Row 45
$ws.Range("D45").Value = '''0.008616'
$ws.Range("E45").Value = '''0.14%'

# Row 46
$ws.Range("D46").Value = '''0.00006578'
$ws.Range("E46").Value = '''-1.64%'

# Row 47
$ws.Range("D47").Value = '''0.00000000752'
$ws.Range("E47").Value = '''0.26%'

# Row 48
$ws.Range("E48").Value = '''13.82%'

# Row 49
$ws.Range("D49").Value = '''0.001446'
$ws.Range("E49").Value = '''20.56%'

# Row 50
$ws.Range("D50").Value = '''0.00002104'
$ws.Range("E50").Value = '''0.26%'

# Row 51
$ws.Range("D51").Value = '''0.0002004'
$ws.Range("E51").Value = '''0.26%'
